$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header names: "_old" -> "_FV2210", "_new" -> "_FV2304"
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel Table ("Table1") with an AutoFilter
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (select row 2 first so Excel freezes above/left of it,
# matching the "freeze top row" behaviour -> pane state "frozen", not "frozenSplit")
[void]($ws.Range("A2").Select())
$ws.Application.ActiveWindow.FreezePanes = $true

# Restore the active cell to A1, matching the original file's selection
[void]($ws.Range("A1").Select())
